$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh - GitHub Actions scheduled update

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.261.93"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.038.38"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.00"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.94"
$ws.Range("E6").Value = "  +2.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.702"
$ws.Range("E7").Value = "  +12.36%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000329"
$ws.Range("E11").Value = "  -3.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.01"
$ws.Range("E12").Value = "  +2.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.685.85"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.66"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.039.54"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.60"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.218.17"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.11"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "98.08"
$ws.Range("E22").Value = "  +3.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("E24").Value = "  +4.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.26"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.16"
$ws.Range("E26").Value = "  -9.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("E27").Value = "  -3.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.82"
$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.84"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.59"
$ws.Range("E30").Value = "  +22.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.41"
$ws.Range("E31").Value = "  -1.75%  "

$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "677.02"
$ws.Range("E33").Value = "  -3.27%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.10"
$ws.Range("E34").Value = "  +2.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.69"
$ws.Range("E35").Value = "  +9.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "66.10"
$ws.Range("E36").Value = "  -2.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.443"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0839"
$ws.Range("E38").Value = "  -7.67%  "

$ws.Range("E39").Value = "  -2.66%  "

$ws.Range("E40").Value = "  -4.60%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("E45").Value = "  +3.36%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.46"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.74"
$ws.Range("E47").Value = "  +6.26%  "

$ws.Range("E48").Value = "  -6.38%  "

$ws.Range("E49").Value = "  -5.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000272"
$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.34"
$ws.Range("E51").Value = "  +1.30%  "
